$wb = $excel.ActiveWorkbook

# Rename the worksheets to reflect the new mantenedor (tipo_deudor)
$wb.Worksheets.Item(1).Name = "tipo_deudor"
$wb.Worksheets.Item(2).Name = "tipo_deudor_rel"
